# Restore C10 on the "Rules" sheet from 18 to 1 (numeric value),
# matching revision 163f97b0fa640bbf8735b9707b9e8025f7e7a236.TEST.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
